$d = $word.ActiveDocument

# Regenerate contracts 18-24: shift every [[PERSON_n]] token down by one
# (n -> n-1) for n = 5..58, so that the duplicate PERSON_5 collapses into
# PERSON_4 and every subsequent identifier closes the resulting gap.
# Processed in ascending order so a freshly written [[PERSON_(n-1)]]
# token is never revisited by a later (higher n) rule.
for ($n = 5; $n -le 58; $n++) {
    $old = "[[PERSON_$n]]"
    $new = "[[PERSON_" + ($n - 1) + "]]"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}
